$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New standings table (row -> Time, Pontos, Jogos, Vitorias, Empates, Derrotas, Saldo de Gols)
$data = @{
    2  = @("Palmeiras",       16, 7, 5, 1, 1, 5)
    3  = @("Flamengo",        14, 7, 4, 2, 1, 12)
    5  = @("Cruzeiro",        13, 7, 4, 1, 2, 2)
    6  = @("Fluminense",      13, 7, 4, 1, 2, 1)
    7  = @("Bahia",           12, 7, 3, 3, 1, 0)
    8  = @("Ceará",           11, 7, 3, 2, 2, 2)
    9  = @("Corinthians",     10, 7, 3, 1, 3, -2)
    10 = @("Internacional",    9, 7, 2, 3, 2, 2)
    11 = @("São Paulo",        9, 7, 1, 6, 0, 1)
    12 = @("Botafogo",         8, 7, 2, 2, 3, 1)
    13 = @("Grêmio",           8, 7, 2, 2, 3, -5)
    14 = @("Vasco da Gama",    7, 7, 2, 1, 4, -3)
    15 = @("Juventude",        7, 6, 2, 1, 3, -7)
    16 = @("Mirassol",         7, 6, 1, 4, 1, 2)
    17 = @("Fortaleza",        7, 7, 1, 4, 2, 0)
    19 = @("Vitória",          6, 7, 1, 3, 3, -3)
    20 = @("Santos",           4, 7, 1, 1, 5, -3)
    21 = @("Sport",             2, 7, 0, 2, 5, -6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}
